{"js": "const body = context.document.body;\n\n// Update the date paragraph (first paragraph in the document body).\nconst paras = body.paragraphs;\nparas.load(\"items\");\n\n// Locate the single table that holds the math problems.\nconst tables = body.tables;\ntables.load(\"items\");\n\nawait context.sync();\n\nconst dateParagraph = paras.items[0];\ndateParagraph.getRange().insertText(\"2025-07-05 Saturday\", Word.InsertLocation.replace);\n\nconst table = tables.items[0];\nconst newValues = [[\"77+14=\", \"97-90=\", \"15+74=\", \"37+28=\", \"54+18=\"], [\"17+6=\", \"56-27=\", \"70+18=\", \"23-15=\", \"18+33=\"], [\"23+28=\", \"1+77=\", \"10+71=\", \"43-43=\", \"98-87=\"], [\"94-88=\", \"47+21=\", \"31+32=\", \"25+31=\", \"63-53=\"], [\"69+8=\", \"13-1=\", \"71-15=\", \"9+12=\", \"75-22=\"], [\"87-48=\", \"17+12=\", \"88-83=\", \"7+16=\", \"81-40=\"], [\"30+8=\", \"99-46=\", \"77-30=\", \"5+45=\", \"52-44=\"], [\"23+61=\", \"85+4=\", \"2+51=\", \"79-21=\", \"4+24=\"], [\"56-2=\", \"33+47=\", \"28-15=\", \"92-32=\", \"57-3=\"], [\"22-7=\", \"99-15=\", \"59-45=\", \"92-33=\", \"60+7=\"], [\"92-62=\", \"91-25=\", \"74-68=\", \"73-19=\", \"33+16=\"], [\"68-2=\", \"52+41=\", \"78-49=\", \"7+77=\", \"14+27=\"], [\"57-0=\", \"69-23=\", \"60-60=\", \"31+55=\", \"36-8=\"], [\"8+8=\", \"95-44=\", \"40+46=\", \"60-14=\", \"16+81=\"], [\"23+54=\", \"50+39=\", \"31+15=\", \"82-6=\", \"74+2=\"], [\"51-11=\", \"62-11=\", \"63-56=\", \"70+29=\", \"48+27=\"], [\"51+47=\", \"36+7=\", \"55-21=\", \"72-19=\", \"20+50=\"], [\"16+12=\", \"47-1=\", \"2+83=\", \"6+55=\", \"36-29=\"], [\"75-7=\", \"8+6=\", \"2+91=\", \"11+0=\", \"48-35=\"], [\"28+0=\", \"91+8=\", \"23+61=\", \"57-1=\", \"30+61=\"]];\n\n// Batch-load the paragraph of every cell first so we only need one more sync.\nconst cellParagraphs = [];\nfor (let r = 0; r < newValues.length; r++) {\n  const rowParagraphs = [];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParas = cell.body.paragraphs;\n    cellParas.load(\"items\");\n    rowParagraphs.push(cellParas);\n  }\n  cellParagraphs.push(rowParagraphs);\n}\nawait context.sync();\n\n// Replace each cell's text in place so run formatting (font, size, alignment) is preserved.\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cellParagraph = cellParagraphs[r][c].items[0];\n    cellParagraph.getRange().insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-05 Saturday\"\n\n# Update every math-problem cell in the table, preserving cell formatting.\n$newValues = @(\n    @('77+14=', '97-90=', '15+74=', '37+28=', '54+18='),\n    @('17+6=', '56-27=', '70+18=', '23-15=', '18+33='),\n    @('23+28=', '1+77=', '10+71=', '43-43=', '98-87='),\n    @('94-88=', '47+21=', '31+32=', '25+31=', '63-53='),\n    @('69+8=', '13-1=', '71-15=', '9+12=', '75-22='),\n    @('87-48=', '17+12=', '88-83=', '7+16=', '81-40='),\n    @('30+8=', '99-46=', '77-30=', '5+45=', '52-44='),\n    @('23+61=', '85+4=', '2+51=', '79-21=', '4+24='),\n    @('56-2=', '33+47=', '28-15=', '92-32=', '57-3='),\n    @('22-7=', '99-15=', '59-45=', '92-33=', '60+7='),\n    @('92-62=', '91-25=', '74-68=', '73-19=', '33+16='),\n    @('68-2=', '52+41=', '78-49=', '7+77=', '14+27='),\n    @('57-0=', '69-23=', '60-60=', '31+55=', '36-8='),\n    @('8+8=', '95-44=', '40+46=', '60-14=', '16+81='),\n    @('23+54=', '50+39=', '31+15=', '82-6=', '74+2='),\n    @('51-11=', '62-11=', '63-56=', '70+29=', '48+27='),\n    @('51+47=', '36+7=', '55-21=', '72-19=', '20+50='),\n    @('16+12=', '47-1=', '2+83=', '6+55=', '36-29='),\n    @('75-7=', '8+6=', '2+91=', '11+0=', '48-35='),\n    @('28+0=', '91+8=', '23+61=', '57-1=', '30+61=')\n)\n\n$tbl = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n  for ($c = 0; $c -lt $newValues[$r].Length; $c++) {\n    $tbl.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n  }\n}\n"}
